# AGA206 Assessment 2 Checklist - mark "Free Camera Mode" and the
# new "Player Camera Controller" task as Done.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Task #22 (row 35) used to be the generic "Come up with your own" task;
# rename it to the actual feature that was implemented.
$ws.Range("C35").Value = "Player Camera Controller "

# Tick the two checkboxes (task #19 "Free Camera Mode" in row 32 and the
# renamed task #22 in row 35) - their linked cells are J32 / J35.
$ws.Range("J32").Value = $true
$ws.Range("J35").Value = $true

# Also flip the underlying form-control state so the checkboxes themselves
# render as checked (in addition to the linked-cell value above).
foreach ($pair in @(@("Check Box 48","J32"), @("Check Box 52","J35"))) {
    $name = $pair[0]
    try {
        $shp = $ws.Shapes.Item($name)
        $shp.ControlFormat.Checked = 1
    } catch {
    }
}
